$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1934.2273
$ws.Range("J17").Value = 1893.0476
$ws.Range("L17").Value = 5679.142800000001
$ws.Range("N17").Value = -6015.142800000001
$ws.Range("H62").Value = 6145.7144
$ws.Range("I62").Value = 6145.7144
$ws.Range("K62").Value = 6145.7144
$ws.Range("M62").Value = -5521.7144
$ws.Range("H65").Value = 6145.7144
$ws.Range("I65").Value = 6145.7144
$ws.Range("K65").Value = 30728.572
$ws.Range("M65").Value = -27608.572
$ws.Range("H86").Value = 90911050
$ws.Range("I86").Value = 166668380
$ws.Range("K86").Value = 166668380
$ws.Range("M86").Value = -166667257
$ws.Range("H89").Value = 90911050
$ws.Range("I89").Value = 166668380
$ws.Range("K89").Value = 833341900
$ws.Range("M89").Value = -833336284
$ws.Range("H96").Value = 257.16666
$ws.Range("I96").Value = 247.33333
$ws.Range("J96").Value = 286.66666
$ws.Range("K96").Value = 741.99999
$ws.Range("L96").Value = 859.9999799999999
$ws.Range("M96").Value = 631.00001
$ws.Range("N96").Value = -3605.99998
$ws.Range("H132").Value = 10450.4
$ws.Range("I132").Value = 11161.214
$ws.Range("J132").Value = 499
$ws.Range("K132").Value = 33483.642
$ws.Range("L132").Value = 1497
$ws.Range("M132").Value = -30953.642
$ws.Range("N132").Value = -6557
$ws.Range("H138").Value = 367138.03
$ws.Range("I138").Value = 6742.5454
$ws.Range("J138").Value = 435488.9
$ws.Range("K138").Value = 20227.6362
$ws.Range("L138").Value = 1306466.7
$ws.Range("M138").Value = -15087.6362
$ws.Range("N138").Value = -1316746.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3131.6667
$ws.Range("I61").Value = 1956.2916
$ws.Range("K61").Value = 1956.2916
$ws.Range("M61").Value = -1744.2916
$ws.Range("H74").Value = 267861.72
$ws.Range("I74").Value = 329916.66
$ws.Range("K74").Value = 329916.66
$ws.Range("M74").Value = -329042.66
$ws.Range("H77").Value = 267861.72
$ws.Range("I77").Value = 329916.66
$ws.Range("K77").Value = 1649583.3
$ws.Range("M77").Value = -1645215.3
$ws.Range("H110").Value = 2648
$ws.Range("I110").Value = 1515.0769
$ws.Range("J110").Value = 4120.8
$ws.Range("K110").Value = 1515.0769
$ws.Range("L110").Value = 4120.8
$ws.Range("M110").Value = 529.9231
$ws.Range("N110").Value = -8210.799999999999
$ws.Range("H122").Value = 2216.875
$ws.Range("I122").Value = 2216.875
$ws.Range("K122").Value = 6650.625
$ws.Range("M122").Value = -4200.625
$ws.Range("H136").Value = 3131.6667
$ws.Range("I136").Value = 1956.2916
$ws.Range("K136").Value = 5868.8748
$ws.Range("M136").Value = -3318.8748

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2790.3142
$ws.Range("I134").Value = 2225.4614
$ws.Range("K134").Value = 6676.3842
$ws.Range("M134").Value = -4141.3842

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1835.4166
$ws.Range("I16").Value = 2032.7
$ws.Range("J16").Value = 849
$ws.Range("K16").Value = 2032.7
$ws.Range("L16").Value = 849
$ws.Range("M16").Value = -1745.7
$ws.Range("N16").Value = -1423
$ws.Range("H31").Value = 3976.5469
$ws.Range("I31").Value = 3343.8918
$ws.Range("K31").Value = 3343.8918
$ws.Range("M31").Value = -3048.8918
$ws.Range("H34").Value = 3976.5469
$ws.Range("I34").Value = 3343.8918
$ws.Range("K34").Value = 3343.8918
$ws.Range("M34").Value = -3141.8918
$ws.Range("H52").Value = 95530.5
$ws.Range("J52").Value = 95530.5
$ws.Range("L52").Value = 95530.5
$ws.Range("N52").Value = -96118.5
$ws.Range("H58").Value = 2947.8215
$ws.Range("I58").Value = 1696.2941
$ws.Range("K58").Value = 1696.2941
$ws.Range("M58").Value = -1493.2941
$ws.Range("H113").Value = 1835.4166
$ws.Range("I113").Value = 2032.7
$ws.Range("J113").Value = 849
$ws.Range("K113").Value = 2032.7
$ws.Range("L113").Value = 849
$ws.Range("M113").Value = 137.3
$ws.Range("N113").Value = -5189
$ws.Range("H132").Value = 4561.0835
$ws.Range("I132").Value = 4077.158
$ws.Range("J132").Value = 6400
$ws.Range("K132").Value = 12231.474
$ws.Range("L132").Value = 19200
$ws.Range("M132").Value = -9701.474
$ws.Range("N132").Value = -24260
$ws.Range("H134").Value = 2390.6052
$ws.Range("I134").Value = 1976.4667
$ws.Range("K134").Value = 5929.4001
$ws.Range("M134").Value = -3394.4001
$ws.Range("H136").Value = 2947.8215
$ws.Range("I136").Value = 1696.2941
$ws.Range("K136").Value = 5088.8823
$ws.Range("M136").Value = -2538.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 999.5
$ws.Range("J107").Value = 1099.7
$ws.Range("L107").Value = 3299.1
$ws.Range("N107").Value = -7139.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H113").Value = 4987.64
$ws.Range("I113").Value = 5193.8945
$ws.Range("K113").Value = 5193.8945
$ws.Range("M113").Value = -3023.8945
$ws.Range("H132").Value = 3822.558
$ws.Range("I132").Value = 3165.889
$ws.Range("J132").Value = 4930.6875
$ws.Range("K132").Value = 9497.667000000001
$ws.Range("L132").Value = 14792.0625
$ws.Range("M132").Value = -6967.667000000001
$ws.Range("N132").Value = -19852.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3698.7
$ws.Range("I61").Value = 4142
$ws.Range("K61").Value = 4142
$ws.Range("M61").Value = -3940
$ws.Range("H113").Value = 3698.7
$ws.Range("I113").Value = 4142
$ws.Range("K113").Value = 4142
$ws.Range("M113").Value = -1972
$ws.Range("H132").Value = 3644.4856
$ws.Range("I132").Value = 2789.0952
$ws.Range("J132").Value = 4927.5713
$ws.Range("K132").Value = 8367.285600000001
$ws.Range("L132").Value = 14782.7139
$ws.Range("M132").Value = -5837.285600000001
$ws.Range("N132").Value = -19842.7139
$ws.Range("H136").Value = 5779.593
$ws.Range("I136").Value = 4503.0625
$ws.Range("K136").Value = 13509.1875
$ws.Range("M136").Value = -10959.1875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26802
$ws.Range("H132").Value = 4686.4595
$ws.Range("I132").Value = 4418.8125
$ws.Range("K132").Value = 13256.4375
$ws.Range("M132").Value = -10726.4375
$ws.Range("H136").Value = 31251320
$ws.Range("I136").Value = 37038330
$ws.Range("J136").Value = 1478
$ws.Range("K136").Value = 111114990
$ws.Range("L136").Value = 4434
$ws.Range("M136").Value = -111112440
$ws.Range("N136").Value = -9534
